$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of locations (night of the big wind additions) appended after
# the existing data (which ends at row 60).
$rows = @(
    @{ Row = 61; Name = "Berlin";        Lat = 52.52;   Lon = 13.4;   NDay = 3; Type = 3 },
    @{ Row = 62; Name = "Herisau";       Lat = 47.385;  Lon = 9.279;  NDay = 3; Type = 1 },
    @{ Row = 63; Name = "Bern";          Lat = 46.948;  Lon = 7.452;  NDay = 2; Type = 1 },
    @{ Row = 64; Name = "Bern Trechsel"; Lat = 46.947;  Lon = 7.451;  NDay = 4; Type = 1 },
    @{ Row = 65; Name = "Fribourg";      Lat = 46.807;  Lon = 7.158;  NDay = 1; Type = 1 },
    @{ Row = 66; Name = "Nufenen";       Lat = 46.5395; Lon = 9.244;  NDay = 3; Type = 1 },
    @{ Row = 67; Name = "St Gallen";     Lat = 47.424;  Lon = 9.378;  NDay = 4; Type = 1 },
    @{ Row = 68; Name = "Schaffhausen";  Lat = 47.696;  Lon = 8.639;  NDay = 3; Type = 1 },
    @{ Row = 69; Name = "Vevey";         Lat = 46.46;   Lon = 6.84;   NDay = 2; Type = 1 },
    @{ Row = 70; Name = "Zurich";        Lat = 47.37;   Lon = 8.474;  NDay = 4; Type = 1 },
    @{ Row = 71; Name = "Basel";         Lat = 47.55;   Lon = 7.591;  NDay = 5; Type = 1 },
    @{ Row = 72; Name = "Rovereto";      Lat = 45.883;  Lon = 11.05;  NDay = 1; Type = 1 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Name
    $ws.Cells.Item($row, 2).Value = $r.Lat
    $ws.Cells.Item($row, 3).Value = $r.Lon
    $ws.Cells.Item($row, 4).Value = $r.NDay
    $ws.Cells.Item($row, 5).Value = $r.Type
}

# Freeze the header row and scroll/select like the saved view in the
# target workbook.
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 55
$null = $ws.Range("C65").Select()
